# Add duplicate detection for contract note imports
#
# A newer trade (2026-02-10, BSE, CN#252611730667) is inserted above the
# existing trade row (2026-02-09, NSE, CN#252611665409), pushing the
# existing trade down from row 5 to row 6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trading History")

# Shift the existing row 5 entry down to row 6 (read with Value2 - Value
# on this host surfaces the property descriptor rather than the scalar).
$ws.Range("A6").Value = $ws.Range("A5").Value2
$ws.Range("B6").Value = $ws.Range("B5").Value2
$ws.Range("C6").Value = $ws.Range("C5").Value2
$ws.Range("D6").Value = $ws.Range("D5").Value2
$ws.Range("E6").Value = $ws.Range("E5").Value2
$ws.Range("F6").Value = $ws.Range("F5").Value2
$ws.Range("G6").Value = $ws.Range("G5").Value2
$ws.Range("H6").Value = $ws.Range("H5").Value2
$ws.Range("I6").Value = $ws.Range("I5").Value2
$ws.Range("J6").Formula = $ws.Range("J5").Formula
$ws.Range("A6").NumberFormat = $ws.Range("A5").NumberFormat

# Write the new (latest) contract-note entry into row 5.
$ws.Range("A5").Value = (Get-Date -Year 2026 -Month 2 -Day 10 -Hour 0 -Minute 0 -Second 0)
$ws.Range("B5").Value = "BSE"
$ws.Range("C5").Value = "Buy"
$ws.Range("D5").Value = 10
$ws.Range("E5").Value = 313.1
$ws.Range("F5").Value = 3153.26
$ws.Range("G5").Value = "CN#252611730667"
$ws.Range("H5").Value = 3.13
$ws.Range("I5").Value = 19.13
$ws.Range("J5").Formula = "=Index!`$C`$2"
